# Insert a new data row before row 149, pushing the existing rows 149..243
# down to 150..244 (dimension grows from A1:R243 to A1:R244), then populate
# the newly inserted row 149 with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 149; everything below shifts down one row.
$ws.Rows.Item(149).Insert()

# The row that used to be 149 is now row 150. Clone its full contents
# (values, shared static columns, styles) into the freshly inserted row 149
# so every column/format matches the rest of the table.
$ws.Rows.Item(150).Copy()
$ws.Rows.Item(149).PasteSpecial()

# Overwrite the columns that differ for this new record: Fecha, Volumen,
# Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg.
$ws.Cells.Item(149, 4).Value = 44777
$ws.Cells.Item(149, 10).Value = 120
$ws.Cells.Item(149, 11).Value = 4000
$ws.Cells.Item(149, 12).Value = 4000
$ws.Cells.Item(149, 13).Value = 4000
$ws.Cells.Item(149, 16).Value = 1333
